# "Generate Report for Handoff"
#
# b.md has finished handoff (for both the zh-cn and de-de locales), so the
# localization-status report needs to move b.md from
# "Handed back: in sync with en-US" to "Ready for handoff", and record the
# new handoff file name + handoff timestamp for each locale.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row 3 is b.md
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"            # zh-cn status
$overview.Range("C3").Value = "Ready for handoff"             # de-de status
$overview.Range("D3").Value = "2016-28-21 00:28:53"           # Latest Handoff Date

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is b.md
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"                 # Status

# D3 carries a hyperlink to the handoff file - update the link's display
# text (this also updates the underlying cell text) so it points at the
# newly generated handoff package instead of the old one.
$zhcnD3Link = $zhcn.Hyperlinks.Item(8)
$zhcnD3Link.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"

$zhcn.Range("E3").Value = "2016-03-21 00:28:50"               # Latest Handoff Datetime

# ---------------------------------------------------------------------
# de-de sheet: row 3 is b.md
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"                 # Status

$dedeD3Link = $dede.Hyperlinks.Item(8)
$dedeD3Link.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"

$dede.Range("E3").Value = "2016-03-21 00:28:53"               # Latest Handoff Datetime
